$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.370833
$ws.Range("H2").Value = 13.112499
$ws.Range("I2").Value = 0.0996525630698175
$ws.Range("J2").Value = 0.09965256306981748
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.716522666666667
$ws.Range("N2").Value = 14.149568
$ws.Range("O2").Value = 0.530211572117814
$ws.Range("P2").Value = 0.530211572117814
$ws.Range("Q2").Value = 20.61513291671467
$ws.Range("R2").Value = 185.536196250432
$ws.Range("S2").Value = 0.05283694213081755
$ws.Range("T2").Value = 0.05283694213081754

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.370833
$ws.Range("H3").Value = 13.112499
$ws.Range("I3").Value = 0.0996525630698175
$ws.Range("J3").Value = 0.09965256306981748
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.236474
$ws.Range("N3").Value = 0.709422
$ws.Range("O3").Value = 0.02658340904223816
$ws.Range("P3").Value = 0.02658340904223817
$ws.Range("Q3").Value = 1.033588362842
$ws.Range("R3").Value = 9.302295265578
$ws.Range("S3").Value = 0.002649104846192395
$ws.Range("T3").Value = 0.002649104846192395

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.370833
$ws.Range("H4").Value = 13.112499
$ws.Range("I4").Value = 0.0996525630698175
$ws.Range("J4").Value = 0.09965256306981748
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.942551666666667
$ws.Range("N4").Value = 11.827655
$ws.Range("O4").Value = 0.4432050188399478
$ws.Range("P4").Value = 0.4432050188399479
$ws.Range("Q4").Value = 17.23223492887167
$ws.Range("R4").Value = 155.090114359845
$ws.Range("S4").Value = 0.04416651609280755
$ws.Range("T4").Value = 0.04416651609280756

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 37.55834333333333
$ws.Range("H5").Value = 112.67503
$ws.Range("I5").Value = 0.8563093528905953
$ws.Range("J5").Value = 0.8563093528905952
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.716522666666667
$ws.Range("N5").Value = 14.149568
$ws.Range("O5").Value = 0.530211572117814
$ws.Range("P5").Value = 0.530211572117814
$ws.Range("Q5").Value = 177.1447776541156
$ws.Range("R5").Value = 1594.30299888704
$ws.Range("S5").Value = 0.4540251282153105
$ws.Range("T5").Value = 0.4540251282153104

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 37.55834333333333
$ws.Range("H6").Value = 112.67503
$ws.Range("I6").Value = 0.8563093528905953
$ws.Range("J6").Value = 0.8563093528905952
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.236474
$ws.Range("N6").Value = 0.709422
$ws.Range("O6").Value = 0.02658340904223816
$ws.Range("P6").Value = 0.02658340904223817
$ws.Range("Q6").Value = 8.881571681406665
$ws.Range("R6").Value = 79.93414513265999
$ws.Range("S6").Value = 0.02276362179458496
$ws.Range("T6").Value = 0.02276362179458497

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 37.55834333333333
$ws.Range("H7").Value = 112.67503
$ws.Range("I7").Value = 0.8563093528905953
$ws.Range("J7").Value = 0.8563093528905952
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.942551666666667
$ws.Range("N7").Value = 11.827655
$ws.Range("O7").Value = 0.4432050188399478
$ws.Range("P7").Value = 0.4432050188399479
$ws.Range("Q7").Value = 148.0757091060722
$ws.Range("R7").Value = 1332.68138195465
$ws.Range("S7").Value = 0.3795206028806998
$ws.Range("T7").Value = 0.3795206028806998

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cadm3"
$ws.Range("C8").Value = "Cadm1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.931542
$ws.Range("H8").Value = 5.794626
$ws.Range("I8").Value = 0.04403808403958729
$ws.Range("J8").Value = 0.04403808403958728
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.716522666666667
$ws.Range("N8").Value = 14.149568
$ws.Range("O8").Value = 0.530211572117814
$ws.Range("P8").Value = 0.530211572117814
$ws.Range("Q8").Value = 9.110161624618668
$ws.Range("R8").Value = 81.991454621568
$ws.Range("S8").Value = 0.02334950177168599
$ws.Range("T8").Value = 0.02334950177168599

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cadm3"
$ws.Range("C9").Value = "Cadm1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.931542
$ws.Range("H9").Value = 5.794626
$ws.Range("I9").Value = 0.04403808403958729
$ws.Range("J9").Value = 0.04403808403958728
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.236474
$ws.Range("N9").Value = 0.709422
$ws.Range("O9").Value = 0.02658340904223816
$ws.Range("P9").Value = 0.02658340904223817
$ws.Range("Q9").Value = 0.456759462908
$ws.Range("R9").Value = 4.110835166172
$ws.Range("S9").Value = 0.001170682401460809
$ws.Range("T9").Value = 0.001170682401460809

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cadm3"
$ws.Range("C10").Value = "Cadm1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.931542
$ws.Range("H10").Value = 5.794626
$ws.Range("I10").Value = 0.04403808403958729
$ws.Range("J10").Value = 0.04403808403958728
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.942551666666667
$ws.Range("N10").Value = 11.827655
$ws.Range("O10").Value = 0.4432050188399478
$ws.Range("P10").Value = 0.4432050188399479
$ws.Range("Q10").Value = 7.615204131336667
$ws.Range("R10").Value = 68.53683718203
$ws.Range("S10").Value = 0.01951789986644049
$ws.Range("T10").Value = 0.01951789986644049

